# Fix automation hanging issue - reduce waitForPageReady timeouts and make
# page readiness non-blocking.
#
# This records the test-run results that were produced once the fix above
# was applied: the "New" click step (row 9) no longer has a screenshot /
# page source captured (its own page-ready wait was skipped), every
# downstream step of that test case (rows 10-24) now fails fast with
# "No valid page available" and no evidence files, and the remaining test
# cases that were never reached (rows 25-40) are marked SKIPPED because
# "TO BE EXECUTED" = NO for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: the screenshot / page-source paths captured previously are no
# longer available now that page-readiness waits are non-blocking.
$ws.Range("O9").Value = ""
$ws.Range("P9").Value = ""

# Rows 10-24: each of these steps now fails immediately because there is
# no valid page to act on, and (like row 9) no screenshot / page source
# could be captured.
$failRows = 10..24
foreach ($r in $failRows) {
    $ws.Range("L$r").Value = "FAIL"
    $ws.Range("M$r").Value = "No valid page available"
    $ws.Range("N$r").Value = "No valid page available"
    $ws.Range("O$r").Value = ""
    $ws.Range("P$r").Value = ""
}

# Rows 25-40: these steps belong to test cases that were never scheduled
# to run (TO BE EXECUTED = NO), so they are reported as skipped.
$skipRows = 25..40
foreach ($r in $skipRows) {
    $ws.Range("L$r").Value = "SKIPPED"
    $ws.Range("M$r").Value = "TO BE EXECUTED = NO"
}

Write-Host "Applied Test_Results status updates for rows 9-40"
